$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column width adjustments (matches target col widths in the diff)
$ws.Columns.Item(1).ColumnWidth = 11.333333333333332
$ws.Columns.Item(2).ColumnWidth = 12.166666666666664
$ws.Columns.Item(3).ColumnWidth = 14.166666666666664
$ws.Columns.Item(4).ColumnWidth = 24.666666666666664

# New "Column7" (G) numeric data for rows 3-102
$ws.Range("G3").Value = 50.73
$ws.Range("G4").Value = 57.122
$ws.Range("G5").Value = 46.984999999999999
$ws.Range("G6").Value = 68.903999999999996
$ws.Range("G7").Value = 35.348999999999997
$ws.Range("G8").Value = 32.594000000000001
$ws.Range("G9").Value = 38.79
$ws.Range("G10").Value = 47.353999999999999
$ws.Range("G11").Value = 42.637
$ws.Range("G12").Value = 44.158999999999999
$ws.Range("G13").Value = 40.957999999999998
$ws.Range("G14").Value = 40.505000000000003
$ws.Range("G15").Value = 34.71
$ws.Range("G16").Value = 44.46
$ws.Range("G17").Value = 75.510000000000005
$ws.Range("G18").Value = 37.774999999999999
$ws.Range("G19").Value = 39.883000000000003
$ws.Range("G20").Value = 53.984000000000002
$ws.Range("G21").Value = 43.45
$ws.Range("G22").Value = 51.164000000000001
$ws.Range("G23").Value = 44.768999999999998
$ws.Range("G24").Value = 55.622
$ws.Range("G25").Value = 55.427999999999997
$ws.Range("G26").Value = 46.991
$ws.Range("G27").Value = 32.549999999999997
$ws.Range("G28").Value = 40.301000000000002
$ws.Range("G29").Value = 27.945
$ws.Range("G30").Value = 31.759
$ws.Range("G31").Value = 44.09
$ws.Range("G32").Value = 41.805
$ws.Range("G33").Value = 54.143999999999998
$ws.Range("G34").Value = 42.021999999999998
$ws.Range("G35").Value = 70.105000000000004
$ws.Range("G36").Value = 76.501999999999995
$ws.Range("G37").Value = 58.457000000000001
$ws.Range("G38").Value = 72.137
$ws.Range("G39").Value = 108.938
$ws.Range("G40").Value = 136.44200000000001
$ws.Range("G41").Value = 114.246
$ws.Range("G42").Value = 69.489999999999995
$ws.Range("G43").Value = 157.529
$ws.Range("G44").Value = 150.964
$ws.Range("G45").Value = 163.71
$ws.Range("G46").Value = 250.45099999999999
$ws.Range("G47").Value = 176.70699999999999
$ws.Range("G48").Value = 272.12599999999998
$ws.Range("G49").Value = 159.685
$ws.Range("G50").Value = 262.48700000000002
$ws.Range("G51").Value = 202.03200000000001
$ws.Range("G52").Value = 161.91300000000001
$ws.Range("G53").Value = 389.15
$ws.Range("G54").Value = 287.17
$ws.Range("G55").Value = 172.22200000000001
$ws.Range("G56").Value = 143.72300000000001
$ws.Range("G57").Value = 162.71199999999999
$ws.Range("G58").Value = 147.59100000000001
$ws.Range("G59").Value = 177.483
$ws.Range("G60").Value = 141.84
$ws.Range("G61").Value = 168.27199999999999
$ws.Range("G62").Value = 164.79
$ws.Range("G63").Value = 170.90299999999999
$ws.Range("G64").Value = 194.83199999999999
$ws.Range("G65").Value = 236.77699999999999
$ws.Range("G66").Value = 279.31799999999998
$ws.Range("G67").Value = 234.07400000000001
$ws.Range("G68").Value = 212.64400000000001
$ws.Range("G69").Value = 321.149
$ws.Range("G70").Value = 226.673
$ws.Range("G71").Value = 268.22699999999998
$ws.Range("G72").Value = 250.97900000000001
$ws.Range("G73").Value = 157.04
$ws.Range("G74").Value = 227.3
$ws.Range("G75").Value = 180.68799999999999
$ws.Range("G76").Value = 192.24799999999999
$ws.Range("G77").Value = 252.178
$ws.Range("G78").Value = 245.24299999999999
$ws.Range("G79").Value = 211.89599999999999
$ws.Range("G80").Value = 219.11600000000001
$ws.Range("G81").Value = 126.831
$ws.Range("G82").Value = 242.11199999999999
$ws.Range("G83").Value = 342.27199999999999
$ws.Range("G84").Value = 240.197
$ws.Range("G85").Value = 451.91500000000002
$ws.Range("G86").Value = 309.03800000000001
$ws.Range("G87").Value = 348.46899999999999
$ws.Range("G88").Value = 409.86500000000001
$ws.Range("G89").Value = 368.471
$ws.Range("G90").Value = 332.01600000000002
$ws.Range("G91").Value = 249.69499999999999
$ws.Range("G92").Value = 311.233
$ws.Range("G93").Value = 243.72499999999999
$ws.Range("G94").Value = 335.22300000000001
$ws.Range("G95").Value = 394.68400000000003
$ws.Range("G96").Value = 270.416
$ws.Range("G97").Value = 466.95600000000002
$ws.Range("G98").Value = 399.02499999999998
$ws.Range("G99").Value = 484.96300000000002
$ws.Range("G100").Value = 536.96600000000001
$ws.Range("G101").Value = 311.99700000000001
$ws.Range("G102").Value = 451.31700000000001

# Update the active selection/view (scroll to top, select G1)
$ws.Range("G1").Select()
